# Fourth Commit Excel Sheet
# Rename Sheet1 -> ValidLogin and populate a sample login table
# (Username/password/HomePage header row + admin/manager/actiTIME row).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "ValidLogin"

# Header row
$ws.Range("A1").Value = "Username"
$ws.Range("B1").Value = "password"
$ws.Range("C1").Value = "HomePage"

# Data row
$ws.Range("A2").Value = "admin"
$ws.Range("B2").Value = "manager"
$ws.Range("C2").Value = "actiTIME - Enter Time-Track"

# Widen column C so the HomePage text fits (~23.55 characters wide)
$ws.Columns.Item(3).ColumnWidth = 22.666666666666668

# Leave the selection on A2, matching the saved view state
$ws.Range("A2").Select()
